$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "35×84=" "98×53="
Replace-Text "49×13=" "28×89="
Replace-Text "96×12=" "11×21="
Replace-Text "51×49=" "67×38="
Replace-Text "66×81=" "35×86="
Replace-Text "43×97=" "15×89="
Replace-Text "55×36=" "96×28="
Replace-Text "18×86=" "68×25="
Replace-Text "42×24=" "19×85="
Replace-Text "83×39=" "20×26="
Replace-Text "76×80=" "29×36="
Replace-Text "95×83=" "21×35="
Replace-Text "93×40=" "86×78="
Replace-Text "83×34=" "14×97="
Replace-Text "40×19=" "98×55="
Replace-Text "87×92=" "74×29="
Replace-Text "81×97=" "59×52="
Replace-Text "43×82=" "59×13="
Replace-Text "50×14=" "64×18="
Replace-Text "31×46=" "53×44="
Replace-Text "17×44=" "27×83="
Replace-Text "36×29=" "28×57="
Replace-Text "13×98=" "23×25="
Replace-Text "83×41=" "44×78="
Replace-Text "88×32=" "45×98="
